$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 504; this shifts the existing rows 504:575
# down to 505:576 and extends the sheet dimension to A1:R576.
$ws.Rows.Item(504).Insert()

# Populate the newly inserted row 504 with its data.
$ws.Cells.Item(504, 1).Value = 3
$ws.Cells.Item(504, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(504, 3).Value = "Coquimbo"
$ws.Cells.Item(504, 4).Value = 45127
$ws.Cells.Item(504, 5).Value = 5
$ws.Cells.Item(504, 6).Value = 100112012
$ws.Cells.Item(504, 7).Value = "Espinaca"
$ws.Cells.Item(504, 8).Value = "Sin especificar"
$ws.Cells.Item(504, 9).Value = "Primera"
$ws.Cells.Item(504, 10).Value = 160
$ws.Cells.Item(504, 11).Value = 4000
$ws.Cells.Item(504, 12).Value = 4500
$ws.Cells.Item(504, 13).Value = 4156
$ws.Cells.Item(504, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(504, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(504, 16).Value = 1385
$ws.Cells.Item(504, 17).Value = 3
$ws.Cells.Item(504, 18).Value = "Hortaliza"
